$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @{C = -0.2389449703450778; D = 0.8125827212666379},
    @{C = -1.171268463748432;  D = 0.24963412143771},
    @{C = -2.417190452640857;  D = 0.02115332047505358},
    @{C = -3.225459107604059;  D = 0.002779086965311706},
    @{C = -0.7969171677649532; D = 0.4310292030191214},
    @{C = -2.144727763237979;  D = 0.03920659739076826},
    @{C = -3.035386725396774;  D = 0.004584682738836721},
    @{C = -2.088343028949126;  D = 0.0443281788042873},
    @{C = -2.638478240736874;  D = 0.01247337720791508},
    @{C = -0.9165226581735678; D = 0.36584833230289}
)

$row = 2
foreach ($entry in $values) {
    $ws.Cells.Item($row, 3).Value = $entry.C
    $ws.Cells.Item($row, 4).Value = $entry.D
    $row++
}
